# Applies the "Added towers, setup tower XML foundation" commit:
#  - rebalances several enemy stat cells on Sheet1 (maxHp/maxMoveSpeed/
#    goldValue/numberOfEnemies/roundEndBonus/freq/armorType); the
#    dependent formula columns (K/L/M/N) recalc automatically.
#  - switches the sheet's frozen pane from a row-freeze (split under
#    row 1) to a column-freeze (split right of column A), and updates
#    the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Enemy stat table edits (columns C..I, rows 2-21) -----------------
$ws.Range("D2").Value = 4

$ws.Range("D3").Value = 4
$ws.Range("H3").Value = 18

$ws.Range("H4").Value = 15

$ws.Range("D5").Value = 4

$ws.Range("D7").Value = 3
$ws.Range("H7").Value = 14

$ws.Range("D8").Value = 3
$ws.Range("H8").Value = 16

$ws.Range("D9").Value = 3

$ws.Range("D10").Value = 4

$ws.Range("H11").Value = 17

$ws.Range("C12").Value = 850
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 20
$ws.Range("H12").Value = 14

$ws.Range("C13").Value = 1200
$ws.Range("D13").Value = 4
$ws.Range("F13").Value = 15
$ws.Range("H13").Value = 14

$ws.Range("C14").Value = 1800
$ws.Range("D14").Value = 4
$ws.Range("E14").Value = 7
$ws.Range("F14").Value = 12

$ws.Range("C15").Value = 2200
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = 14
$ws.Range("H15").Value = 18

$ws.Range("C16").Value = 2600
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 7

$ws.Range("C17").Value = 3200
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 7
$ws.Range("H17").Value = 14

$ws.Range("C18").Value = 3600
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 8
$ws.Range("G18").Value = 220
$ws.Range("H18").Value = 15
$ws.Range("I18").Value = "light"

$ws.Range("C19").Value = 4200
$ws.Range("E19").Value = 9
$ws.Range("G19").Value = 240
$ws.Range("I19").Value = "heavy"

$ws.Range("C20").Value = 4800
$ws.Range("E20").Value = 10
$ws.Range("G20").Value = 260
$ws.Range("H20").Value = 15

$ws.Range("C21").Value = 6000
$ws.Range("D21").Value = 4
$ws.Range("E21").Value = 10
$ws.Range("G21").Value = 280

# --- View: swap the frozen-row pane for a frozen-column pane ---------
# Unfreeze first so the split can be redefined, then select the cell
# just right of the column we want frozen (A) and freeze there, which
# yields xSplit=1 / activePane="topRight" (mirrors the target pane's
# column freeze instead of the original row freeze under row 1).
$excel.ActiveWindow.FreezePanes = $false
[void]$ws.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true

# Move/extend the selection to the new active cell in the unfrozen pane.
[void]$ws.Range("F8").Select()
